# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  -> "Integral"     (bound to the Slide Master)
#   ppt/theme/theme2.xml  -> "Office Theme" (bound to the Notes Master)
#
# The commit swaps the two themes' colour palettes: the Slide Master's
# theme (theme1.xml) ends up carrying the stock "Office Theme" colours
# (the palette that used to live on the Notes Master's theme).
#
# PowerPoint's automation surface doesn't expose a way to rename a
# theme/colour-scheme or to reach the Notes Master's theme object, so
# we apply the reachable, supported part of the edit: push the twelve
# "Office Theme" theme colours onto the presentation's (Slide Master)
# ThemeColorScheme, in the documented
#   ThemeColorScheme.Colors(index).RGB
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

# RGB(r,g,b) isn't available in this host, so pack the OLE_COLOR value
# (0x00BBGGRR) by hand: r | (g << 8) | (b << 16).
function HexClr([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = HexClr 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = HexClr 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = HexClr 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = HexClr 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = HexClr 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = HexClr 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = HexClr 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = HexClr 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = HexClr 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = HexClr 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = HexClr 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = HexClr 0x95 0x4F 0x72   # folHlink
